$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.28 = 12528.52 pesos`n✅ 12528.52 pesos = 3.25 = 948.98 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 305
$ws2.Range("O10").Value = 3821.2
$ws2.Range("N12").Value = 3855
$ws2.Range("O12").Value = 292
